# ---------------------------------------------------------------------------
# SB Squares Task Report - add "Live score auto-update from ESPN" task batch
#
# The 12 new task rows (Task # 114-125) land at rows 95-106. Three of those
# rows (95-97) already existed (blank / Summary header / Total Tasks), so we
# only need to INSERT 9 brand-new rows (98-106) to make room; that single
# insert pushes the existing "Completed:" / "By Assignee" / "By Type" summary
# blocks down by 9 rows (98->107 ... 128->137), which is exactly what the
# target workbook shows, with their original values/styles intact.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 whole rows at 98:106 - shifts old rows 98-128 down to 107-137.
$ws.Range("98:106").Insert()

# Fill rows 95-106 with the new ESPN live-score feature tasks.

# Row 95: Task #114
$ws.Cells.Item(95, 1).Value = 114
$ws.Cells.Item(95, 2).Value = 'Add auto_scores_enabled and espn_event_id columns to games table + update TypeScript types'
$ws.Cells.Item(95, 3).Value = 'Migration'
$ws.Cells.Item(95, 4).Value = 'ui-dev'
$ws.Cells.Item(95, 5).Value = 'Completed'
$ws.Cells.Item(95, 6).Value = 'supabase/migrations/20250207000000_live_scores.sql, lib/types.ts'

# Row 96: Task #115
$ws.Cells.Item(96, 1).Value = 115
$ws.Cells.Item(96, 2).Value = 'Build ESPN API client with scoreboard fetch, game parser, team finder, and team-axis matcher'
$ws.Cells.Item(96, 3).Value = 'Feature'
$ws.Cells.Item(96, 4).Value = 'ui-dev'
$ws.Cells.Item(96, 5).Value = 'Completed'
$ws.Cells.Item(96, 6).Value = 'lib/espn.ts'

# Row 97: Task #116
$ws.Cells.Item(97, 1).Value = 116
$ws.Cells.Item(97, 2).Value = 'Create live-scores API route with ESPN polling, score mapping, and Supabase upserts'
$ws.Cells.Item(97, 3).Value = 'Feature'
$ws.Cells.Item(97, 4).Value = 'ui-dev'
$ws.Cells.Item(97, 5).Value = 'Completed'
$ws.Cells.Item(97, 6).Value = 'app/api/live-scores/route.ts'

# Row 98: Task #117
$ws.Cells.Item(98, 1).Value = 117
$ws.Cells.Item(98, 2).Value = 'Create client-side polling hook with adaptive intervals (30s live, 60s idle)'
$ws.Cells.Item(98, 3).Value = 'Feature'
$ws.Cells.Item(98, 4).Value = 'ui-dev'
$ws.Cells.Item(98, 5).Value = 'Completed'
$ws.Cells.Item(98, 6).Value = 'hooks/use-live-scores.ts'

# Row 99: Task #118
$ws.Cells.Item(99, 1).Value = 118
$ws.Cells.Item(99, 2).Value = 'Add auto-update ESPN toggle, NFL status indicator, and disable manual inputs in ScoreEntry'
$ws.Cells.Item(99, 3).Value = 'Feature'
$ws.Cells.Item(99, 4).Value = 'ui-dev'
$ws.Cells.Item(99, 5).Value = 'Completed'
$ws.Cells.Item(99, 6).Value = 'app/game/[gameId]/admin/page.tsx'

# Row 100: Task #119
$ws.Cells.Item(100, 1).Value = 119
$ws.Cells.Item(100, 2).Value = 'Research NFL live score APIs (ESPN public scoreboard, no auth required)'
$ws.Cells.Item(100, 3).Value = 'Research'
$ws.Cells.Item(100, 4).Value = 'product-researcher'
$ws.Cells.Item(100, 5).Value = 'Completed'
$ws.Cells.Item(100, 6).Value = 'N/A (research only)'

# Row 101: Task #120
$ws.Cells.Item(101, 1).Value = 120
$ws.Cells.Item(101, 2).Value = 'Create architecture plan for live score auto-update feature (data flow, file list, edge cases)'
$ws.Cells.Item(101, 3).Value = 'Docs'
$ws.Cells.Item(101, 4).Value = 'architect'
$ws.Cells.Item(101, 5).Value = 'Completed'
$ws.Cells.Item(101, 6).Value = 'N/A (architecture plan)'

# Row 102: Task #121
$ws.Cells.Item(102, 1).Value = 121
$ws.Cells.Item(102, 2).Value = 'Add Super Bowl auto-detect endpoint using ESPN scoreboard date filtering and team matching'
$ws.Cells.Item(102, 3).Value = 'Feature'
$ws.Cells.Item(102, 4).Value = 'architect'
$ws.Cells.Item(102, 5).Value = 'Completed'
$ws.Cells.Item(102, 6).Value = 'app/api/live-scores/detect/route.ts, lib/espn.ts'

# Row 103: Task #122
$ws.Cells.Item(103, 1).Value = 122
$ws.Cells.Item(103, 2).Value = 'Fix ESPN fetchESPNScores to fallback to summary endpoint for historical games'
$ws.Cells.Item(103, 3).Value = 'Bugfix'
$ws.Cells.Item(103, 4).Value = 'team-lead-2'
$ws.Cells.Item(103, 5).Value = 'Completed'
$ws.Cells.Item(103, 6).Value = 'lib/espn.ts'

# Row 104: Task #123
$ws.Cells.Item(104, 1).Value = 123
$ws.Cells.Item(104, 2).Value = 'Free tier safety analysis for live score polling (30s interval confirmed safe)'
$ws.Cells.Item(104, 3).Value = 'Research'
$ws.Cells.Item(104, 4).Value = 'architect'
$ws.Cells.Item(104, 5).Value = 'Completed'
$ws.Cells.Item(104, 6).Value = 'N/A (analysis only)'

# Row 105: Task #124
$ws.Cells.Item(105, 1).Value = 124
$ws.Cells.Item(105, 2).Value = 'Add live score polling free tier analysis to research findings doc'
$ws.Cells.Item(105, 3).Value = 'Docs'
$ws.Cells.Item(105, 4).Value = 'product-researcher'
$ws.Cells.Item(105, 5).Value = 'Completed'
$ws.Cells.Item(105, 6).Value = 'docs/research-findings.md'

# Row 106: Task #125
$ws.Cells.Item(106, 1).Value = 125
$ws.Cells.Item(106, 2).Value = 'Add Session 5 agent team success story to README'
$ws.Cells.Item(106, 3).Value = 'Docs'
$ws.Cells.Item(106, 4).Value = 'ui-dev'
$ws.Cells.Item(106, 5).Value = 'Completed'
$ws.Cells.Item(106, 6).Value = 'README.md'

Write-Output "Inserted 9 rows and populated tasks 114-125 (rows 95-106)."
